$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows where Target cluster (column D) was "Resolving-Mac" (shared string 23).
# Delete from bottom to top so row indices of earlier rows remain stable.
$ws.Rows.Item(17).Delete() | Out-Null
$ws.Rows.Item(13).Delete() | Out-Null
$ws.Rows.Item(9).Delete() | Out-Null
$ws.Rows.Item(5).Delete() | Out-Null

# Update the remaining rows with the recomputed TPM-based values.

# Row 2: ECs -> ECs
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 5.230988666666666
$ws.Cells.Item(2, 8).Value = 15.692966
$ws.Cells.Item(2, 9).Value = 0.2129406655351238
$ws.Cells.Item(2, 10).Value = 0.2129406655351238
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.227228
$ws.Cells.Item(2, 14).Value = 0.681684
$ws.Cells.Item(2, 15).Value = 0.2376267857721762
$ws.Cells.Item(2, 16).Value = 0.2376267857721762
$ws.Cells.Item(2, 17).Value = 1.188627092749333
$ws.Cells.Item(2, 18).Value = 10.697643834744
$ws.Cells.Item(2, 19).Value = 0.05060040591129947
$ws.Cells.Item(2, 20).Value = 0.05060040591129947

# Row 3: ECs -> FAPs
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 5.230988666666666
$ws.Cells.Item(3, 8).Value = 15.692966
$ws.Cells.Item(3, 9).Value = 0.2129406655351238
$ws.Cells.Item(3, 10).Value = 0.2129406655351238
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 0.6537306666666667
$ws.Cells.Item(3, 14).Value = 1.961192
$ws.Cells.Item(3, 15).Value = 0.6836477770376096
$ws.Cells.Item(3, 16).Value = 0.6836477770376095
$ws.Cells.Item(3, 17).Value = 3.419657708385778
$ws.Cells.Item(3, 18).Value = 30.776919375472
$ws.Cells.Item(3, 19).Value = 0.1455764126339965
$ws.Cells.Item(3, 20).Value = 0.1455764126339965

# Row 4: ECs -> MuSCs
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 5.230988666666666
$ws.Cells.Item(4, 8).Value = 15.692966
$ws.Cells.Item(4, 9).Value = 0.2129406655351238
$ws.Cells.Item(4, 10).Value = 0.2129406655351238
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.07528033333333332
$ws.Cells.Item(4, 14).Value = 0.225841
$ws.Cells.Item(4, 15).Value = 0.0787254371902143
$ws.Cells.Item(4, 16).Value = 0.0787254371902143
$ws.Cells.Item(4, 17).Value = 0.3937905704895555
$ws.Cells.Item(4, 18).Value = 3.544115134405999
$ws.Cells.Item(4, 19).Value = 0.01676384698982782
$ws.Cells.Item(4, 20).Value = 0.01676384698982782

# Row 5: FAPs -> ECs
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 8.489274999999999
$ws.Cells.Item(5, 8).Value = 25.467825
$ws.Cells.Item(5, 9).Value = 0.345577477529236
$ws.Cells.Item(5, 10).Value = 0.3455774775292359
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.227228
$ws.Cells.Item(5, 14).Value = 0.681684
$ws.Cells.Item(5, 15).Value = 0.2376267857721762
$ws.Cells.Item(5, 16).Value = 0.2376267857721762
$ws.Cells.Item(5, 17).Value = 1.9290009797
$ws.Cells.Item(5, 18).Value = 17.3610088173
$ws.Cells.Item(5, 19).Value = 0.08211846522052878
$ws.Cells.Item(5, 20).Value = 0.08211846522052876

# Row 6: FAPs -> FAPs
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 8.489274999999999
$ws.Cells.Item(6, 8).Value = 25.467825
$ws.Cells.Item(6, 9).Value = 0.345577477529236
$ws.Cells.Item(6, 10).Value = 0.3455774775292359
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.6537306666666667
$ws.Cells.Item(6, 14).Value = 1.961192
$ws.Cells.Item(6, 15).Value = 0.6836477770376096
$ws.Cells.Item(6, 16).Value = 0.6836477770376095
$ws.Cells.Item(6, 17).Value = 5.549699405266666
$ws.Cells.Item(6, 18).Value = 49.9472946474
$ws.Cells.Item(6, 19).Value = 0.2362532743071266
$ws.Cells.Item(6, 20).Value = 0.2362532743071266

# Row 7: FAPs -> MuSCs
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 4).Value = "MuSCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 8.489274999999999
$ws.Cells.Item(7, 8).Value = 25.467825
$ws.Cells.Item(7, 9).Value = 0.345577477529236
$ws.Cells.Item(7, 10).Value = 0.3455774775292359
$ws.Cells.Item(7, 11).Value = 1
$ws.Cells.Item(7, 12).Value = 0.3333333333333333
$ws.Cells.Item(7, 13).Value = 0.07528033333333332
$ws.Cells.Item(7, 14).Value = 0.225841
$ws.Cells.Item(7, 15).Value = 0.0787254371902143
$ws.Cells.Item(7, 16).Value = 0.0787254371902143
$ws.Cells.Item(7, 17).Value = 0.6390754517583332
$ws.Cells.Item(7, 18).Value = 5.751679065824999
$ws.Cells.Item(7, 19).Value = 0.02720573800158056
$ws.Cells.Item(7, 20).Value = 0.02720573800158056

# Row 8: MuSCs -> ECs
$ws.Cells.Item(8, 1).Value = "MuSCs"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 8.418577333333333
$ws.Cells.Item(8, 8).Value = 25.255732
$ws.Cells.Item(8, 9).Value = 0.3426995496362334
$ws.Cells.Item(8, 10).Value = 0.3426995496362334
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 0.227228
$ws.Cells.Item(8, 14).Value = 0.681684
$ws.Cells.Item(8, 15).Value = 0.2376267857721762
$ws.Cells.Item(8, 16).Value = 0.2376267857721762
$ws.Cells.Item(8, 17).Value = 1.912936490298667
$ws.Cells.Item(8, 18).Value = 17.216428412688
$ws.Cells.Item(8, 19).Value = 0.08143459246563049
$ws.Cells.Item(8, 20).Value = 0.08143459246563048

# Row 9: MuSCs -> FAPs
$ws.Cells.Item(9, 1).Value = "MuSCs"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 8.418577333333333
$ws.Cells.Item(9, 8).Value = 25.255732
$ws.Cells.Item(9, 9).Value = 0.3426995496362334
$ws.Cells.Item(9, 10).Value = 0.3426995496362334
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.6537306666666667
$ws.Cells.Item(9, 14).Value = 1.961192
$ws.Cells.Item(9, 15).Value = 0.6836477770376096
$ws.Cells.Item(9, 16).Value = 0.6836477770376095
$ws.Cells.Item(9, 17).Value = 5.503482172504889
$ws.Cells.Item(9, 18).Value = 49.531339552544
$ws.Cells.Item(9, 19).Value = 0.2342857853006009
$ws.Cells.Item(9, 20).Value = 0.2342857853006009

# Row 10: MuSCs -> MuSCs
$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 4).Value = "MuSCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 8.418577333333333
$ws.Cells.Item(10, 8).Value = 25.255732
$ws.Cells.Item(10, 9).Value = 0.3426995496362334
$ws.Cells.Item(10, 10).Value = 0.3426995496362334
$ws.Cells.Item(10, 11).Value = 1
$ws.Cells.Item(10, 12).Value = 0.3333333333333333
$ws.Cells.Item(10, 13).Value = 0.07528033333333332
$ws.Cells.Item(10, 14).Value = 0.225841
$ws.Cells.Item(10, 15).Value = 0.0787254371902143
$ws.Cells.Item(10, 16).Value = 0.0787254371902143
$ws.Cells.Item(10, 17).Value = 0.6337533078457777
$ws.Cells.Item(10, 18).Value = 5.703779770612
$ws.Cells.Item(10, 19).Value = 0.02697917187000202
$ws.Cells.Item(10, 20).Value = 0.02697917187000202

# Row 11: Resolving-Mac -> ECs
$ws.Cells.Item(11, 1).Value = "Resolving-Mac"
$ws.Cells.Item(11, 4).Value = "ECs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 2.426634333333333
$ws.Cells.Item(11, 8).Value = 7.279902999999999
$ws.Cells.Item(11, 9).Value = 0.09878230729940689
$ws.Cells.Item(11, 10).Value = 0.09878230729940687
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 0.6666666666666666
$ws.Cells.Item(11, 13).Value = 0.227228
$ws.Cells.Item(11, 14).Value = 0.681684
$ws.Cells.Item(11, 15).Value = 0.2376267857721762
$ws.Cells.Item(11, 16).Value = 0.2376267857721762
$ws.Cells.Item(11, 17).Value = 0.5513992662946665
$ws.Cells.Item(11, 18).Value = 4.962593396651999
$ws.Cells.Item(11, 19).Value = 0.02347332217471744
$ws.Cells.Item(11, 20).Value = 0.02347332217471743

# Row 12: Resolving-Mac -> FAPs
$ws.Cells.Item(12, 1).Value = "Resolving-Mac"
$ws.Cells.Item(12, 4).Value = "FAPs"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 2.426634333333333
$ws.Cells.Item(12, 8).Value = 7.279902999999999
$ws.Cells.Item(12, 9).Value = 0.09878230729940689
$ws.Cells.Item(12, 10).Value = 0.09878230729940687
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 0.6537306666666667
$ws.Cells.Item(12, 14).Value = 1.961192
$ws.Cells.Item(12, 15).Value = 0.6836477770376096
$ws.Cells.Item(12, 16).Value = 0.6836477770376095
$ws.Cells.Item(12, 17).Value = 1.586365280486222
$ws.Cells.Item(12, 18).Value = 14.277287524376
$ws.Cells.Item(12, 19).Value = 0.06753230479588555
$ws.Cells.Item(12, 20).Value = 0.06753230479588554

# Row 13: Resolving-Mac -> MuSCs
$ws.Cells.Item(13, 1).Value = "Resolving-Mac"
$ws.Cells.Item(13, 4).Value = "MuSCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 2.426634333333333
$ws.Cells.Item(13, 8).Value = 7.279902999999999
$ws.Cells.Item(13, 9).Value = 0.09878230729940689
$ws.Cells.Item(13, 10).Value = 0.09878230729940687
$ws.Cells.Item(13, 11).Value = 1
$ws.Cells.Item(13, 12).Value = 0.3333333333333333
$ws.Cells.Item(13, 13).Value = 0.07528033333333332
$ws.Cells.Item(13, 14).Value = 0.225841
$ws.Cells.Item(13, 15).Value = 0.0787254371902143
$ws.Cells.Item(13, 16).Value = 0.0787254371902143
$ws.Cells.Item(13, 17).Value = 0.1826778414914444
$ws.Cells.Item(13, 18).Value = 1.644100573423
$ws.Cells.Item(13, 19).Value = 0.007776680328803904
$ws.Cells.Item(13, 20).Value = 0.007776680328803904
